# Registree stats backup on Tue 13 Apr 2021 21:27:29 SAST
# Adds a new registree (Khoza, Faniso) to the MD410 Attendance sheet and
# refreshes the "as of" timestamps on all report sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet: MD410 Attendance
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("MD410 Attendance")

# Refresh the report-generation timestamp shown in the title cell.
$ws1.Range("A1").Value = "MD410 Registrees as of 13/04/2021 21:27"

# Insert a new row for the new registree so the alphabetical order by
# Last Name is preserved (between "Khan" and "Kienast").
$ws1.Rows.Item(91).Insert()

$ws1.Range("A91").Value = "Khoza"
$ws1.Range("B91").Value = "Faniso"
$ws1.Range("C91").Value = "The Wilds"
$ws1.Range("D91").Value = "Yes"
$ws1.Range("E91").Value = "No"
$ws1.Range("F91").Value = "410E"

# Match the formatting used by the rest of the data rows (thin border
# border box and 25pt row height).
$ws1.Range("A91:F91").Borders.LineStyle = 1
$ws1.Rows.Item(91).RowHeight = 25

# Bump the attendee count shown in the trailing summary row (it shifted
# from row 230 to row 231 because of the inserted row above).
$ws1.Range("A231").Value = "Number of attendees: 228"

# ---------------------------------------------------------------
# Sheet: 410E Attendance
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("410E Attendance")
$ws2.Range("A1").Value = "410E Registrees as of 13/04/2021 21:27"

# ---------------------------------------------------------------
# Sheet: 410W Attendance
# ---------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("410W Attendance")
$ws3.Range("A1").Value = "410W Registrees as of 13/04/2021 21:27"

# ---------------------------------------------------------------
# Sheet: 410E Voting
# ---------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("410E Voting")
$ws4.Range("A1").Value = "410E Voting details as of 13/04/2021 21:27"

# ---------------------------------------------------------------
# Sheet: 410W Voting
# ---------------------------------------------------------------
$ws5 = $wb.Worksheets.Item("410W Voting")
$ws5.Range("A1").Value = "410W Voting details as of 13/04/2021 21:27"
